# Arreglo cosas del ejecutable
# Add a new data row (row 2) to Sheet1: ALM / B1 / 15 / +543413654789
# "15" and the phone number must be stored as text (not numbers), and
# without leaving a persistent cell style, so we briefly force a text
# number format, assign the values, then reset the style back to Normal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ALM"
$ws.Range("B2").Value = "B1"

$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C2").Value = "15"
$ws.Range("D2").Value = "+543413654789"
$ws.Range("C2:D2").Style = "Normal"
